# Auto-generated edit script: applies numeric cell updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 720273.5
$ws.Range("I112").Value = 866.6667
$ws.Range("J112").Value = 803282
$ws.Range("K112").Value = 2600.0001
$ws.Range("L112").Value = 2409846
$ws.Range("M112").Value = -1492.0001
$ws.Range("N112").Value = -2412062
$ws.Range("H137").Value = 14216630
$ws.Range("I137").Value = 989.7742
$ws.Range("J137").Value = 48115464
$ws.Range("K137").Value = 2969.3226
$ws.Range("L137").Value = 144346392
$ws.Range("M137").Value = -419.3226
$ws.Range("N137").Value = -144351492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 11033.77
$ws.Range("I37").Value = 3780
$ws.Range("K37").Value = 3780
$ws.Range("M37").Value = -3507
$ws.Range("H55").Value = 14499
$ws.Range("J55").Value = 14499
$ws.Range("L55").Value = 14499
$ws.Range("N55").Value = -15129
$ws.Range("H61").Value = 3242521.8
$ws.Range("I61").Value = 1812493.4
$ws.Range("J61").Value = 7353853.5
$ws.Range("K61").Value = 1812493.4
$ws.Range("L61").Value = 7353853.5
$ws.Range("M61").Value = -1812281.4
$ws.Range("N61").Value = -7354277.5
$ws.Range("H74").Value = 28396420
$ws.Range("I74").Value = 27027754
$ws.Range("K74").Value = 27027754
$ws.Range("M74").Value = -27026880
$ws.Range("H77").Value = 28396420
$ws.Range("I77").Value = 27027754
$ws.Range("K77").Value = 135138770
$ws.Range("M77").Value = -135134402
$ws.Range("H80").Value = 23124.875
$ws.Range("J80").Value = 23124.875
$ws.Range("L80").Value = 23124.875
$ws.Range("N80").Value = -25120.875
$ws.Range("H83").Value = 23124.875
$ws.Range("J83").Value = 23124.875
$ws.Range("L83").Value = 69374.625
$ws.Range("N83").Value = -79358.625
$ws.Range("H136").Value = 3242521.8
$ws.Range("I136").Value = 1812493.4
$ws.Range("J136").Value = 7353853.5
$ws.Range("K136").Value = 5437480.199999999
$ws.Range("L136").Value = 22061560.5
$ws.Range("M136").Value = -5434930.199999999
$ws.Range("N136").Value = -22066660.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14480017
$ws.Range("I134").Value = 18519582
$ws.Range("J134").Value = 3573189.8
$ws.Range("K134").Value = 55558746
$ws.Range("L134").Value = 10719569.4
$ws.Range("M134").Value = -55556211
$ws.Range("N134").Value = -10724639.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1053933
$ws.Range("I22").Value = 1250233
$ws.Range("J22").Value = 7000
$ws.Range("K22").Value = 1250233
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = -1249883
$ws.Range("N22").Value = -7700
$ws.Range("H31").Value = 1361995.9
$ws.Range("I31").Value = 1491.1471
$ws.Range("J31").Value = 5216759.5
$ws.Range("K31").Value = 1491.1471
$ws.Range("L31").Value = 5216759.5
$ws.Range("M31").Value = -1196.1471
$ws.Range("N31").Value = -5217349.5
$ws.Range("H34").Value = 1361995.9
$ws.Range("I34").Value = 1491.1471
$ws.Range("J34").Value = 5216759.5
$ws.Range("K34").Value = 1491.1471
$ws.Range("L34").Value = 5216759.5
$ws.Range("M34").Value = -1289.1471
$ws.Range("N34").Value = -5217163.5
$ws.Range("H58").Value = 723820.1
$ws.Range("I58").Value = 2828.9565
$ws.Range("K58").Value = 2828.9565
$ws.Range("M58").Value = -2625.9565
$ws.Range("H132").Value = 1325.4131
$ws.Range("I132").Value = 1088.5897
$ws.Range("K132").Value = 3265.7691
$ws.Range("M132").Value = -735.7691
$ws.Range("H134").Value = 582517.75
$ws.Range("I134").Value = 2403.246
$ws.Range("J134").Value = 10009378
$ws.Range("K134").Value = 7209.738
$ws.Range("L134").Value = 30028134
$ws.Range("M134").Value = -4674.738
$ws.Range("N134").Value = -30033204
$ws.Range("H136").Value = 723820.1
$ws.Range("I136").Value = 2828.9565
$ws.Range("K136").Value = 8486.869499999999
$ws.Range("M136").Value = -5936.869499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 35716332
$ws.Range("I129").Value = 83334440
$ws.Range("J129").Value = 5955012
$ws.Range("K129").Value = 250003320
$ws.Range("L129").Value = 17865036
$ws.Range("M129").Value = -249998320
$ws.Range("N129").Value = -17875036
$ws.Range("H130").Value = 1520
$ws.Range("I130").Value = 1765
$ws.Range("J130").Value = 1450
$ws.Range("K130").Value = 5295
$ws.Range("L130").Value = 4350
$ws.Range("M130").Value = -275
$ws.Range("N130").Value = -14390
$ws.Range("H131").Value = 7813315.5
$ws.Range("J131").Value = 977.31915
$ws.Range("L131").Value = 2931.95745
$ws.Range("N131").Value = -13011.95745
$ws.Range("H132").Value = 3715.7
$ws.Range("I132").Value = 3860.5
$ws.Range("J132").Value = 3498.5
$ws.Range("K132").Value = 34744.5
$ws.Range("L132").Value = 31486.5
$ws.Range("M132").Value = -32214.5
$ws.Range("N132").Value = -36546.5
$ws.Range("H133").Value = 2987.7083
$ws.Range("I133").Value = 2977.3076
$ws.Range("K133").Value = 8931.9228
$ws.Range("M133").Value = -3871.9228
$ws.Range("H134").Value = 1700
$ws.Range("I134").Value = 1700
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5100
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -30
$ws.Range("H136").Value = 2269.9524
$ws.Range("I136").Value = 1283.8462
$ws.Range("J136").Value = 3872.375
$ws.Range("K136").Value = 3851.5386
$ws.Range("L136").Value = 11617.125
$ws.Range("M136").Value = 1248.4614
$ws.Range("N136").Value = -21817.125
$ws.Range("H137").Value = 2395.375
$ws.Range("I137").Value = 1911.4286
$ws.Range("J137").Value = 3072.9
$ws.Range("K137").Value = 5734.2858
$ws.Range("L137").Value = 9218.700000000001
$ws.Range("M137").Value = -634.2857999999997
$ws.Range("N137").Value = -19418.7
$ws.Range("H138").Value = 125002190
$ws.Range("I138").Value = 166667920
$ws.Range("J138").Value = 5013
$ws.Range("K138").Value = 500003760
$ws.Range("L138").Value = 15039
$ws.Range("M138").Value = -499998620
$ws.Range("N138").Value = -25319
$ws.Range("H139").Value = 86282.164
$ws.Range("I139").Value = 102532
$ws.Range("J139").Value = 5033
$ws.Range("K139").Value = 307596
$ws.Range("L139").Value = 15099
$ws.Range("M139").Value = -302456
$ws.Range("N139").Value = -25379
$ws.Range("H140").Value = 3000.9688
$ws.Range("I140").Value = 2850.7896
$ws.Range("K140").Value = 8552.3688
$ws.Range("M140").Value = -3372.3688
$ws.Range("H141").Value = 2122.3076
$ws.Range("I141").Value = 1780.909
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 5342.727000000001
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -162.7270000000008
$ws.Range("N141").Value = -22360
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3219.8572
$ws.Range("I22").Value = 2212.2
$ws.Range("J22").Value = 3975.6
$ws.Range("K22").Value = 2212.2
$ws.Range("L22").Value = 3975.6
$ws.Range("M22").Value = -1917.2
$ws.Range("N22").Value = -4565.6
$ws.Range("H27").Value = 3219.8572
$ws.Range("I27").Value = 2212.2
$ws.Range("J27").Value = 3975.6
$ws.Range("K27").Value = 2212.2
$ws.Range("L27").Value = 3975.6
$ws.Range("M27").Value = -2105.2
$ws.Range("N27").Value = -4189.6
$ws.Range("H132").Value = 2919572
$ws.Range("I132").Value = 3864783.5
$ws.Range("K132").Value = 11594350.5
$ws.Range("M132").Value = -11591820.5
$ws.Range("H136").Value = 1764512.9
$ws.Range("I136").Value = 1950163.4
$ws.Range("K136").Value = 5850490.199999999
$ws.Range("M136").Value = -5847940.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2523.5615
$ws.Range("I136").Value = 602.7381
$ws.Range("J136").Value = 5125.968
$ws.Range("K136").Value = 1808.2143
$ws.Range("L136").Value = 15377.904
$ws.Range("M136").Value = 741.7856999999999
$ws.Range("N136").Value = -20477.904

